{"js": "// Replace the date heading and every two-digit-by-two-digit multiplication\n// problem/answer in the table with the new values from the target revision.\n// Every \"old\" value below is unique in the document, so a plain\n// search-and-replace (one search per pair) is sufficient and unambiguous.\nconst replacements = [\n  [\"2025-05-28 Wednesday\", \"2025-05-29 Thursday\"],\n  [\"64\u00d784=5376\", \"65\u00d747=3055\"],\n  [\"63\u00d758=3654\", \"45\u00d789=4005\"],\n  [\"91\u00d723=2093\", \"36\u00d796=3456\"],\n  [\"14\u00d755=770\", \"20\u00d728=560\"],\n  [\"63\u00d738=2394\", \"45\u00d763=2835\"],\n  [\"18\u00d713=234\", \"18\u00d773=1314\"],\n  [\"96\u00d782=7872\", \"59\u00d781=4779\"],\n  [\"43\u00d714=602\", \"61\u00d755=3355\"],\n  [\"46\u00d731=1426\", \"71\u00d745=3195\"],\n  [\"93\u00d747=4371\", \"71\u00d779=5609\"],\n  [\"72\u00d798=7056\", \"13\u00d718=234\"],\n  [\"97\u00d788=8536\", \"37\u00d711=407\"],\n  [\"78\u00d768=5304\", \"17\u00d720=340\"],\n  [\"76\u00d747=3572\", \"28\u00d793=2604\"],\n  [\"17\u00d763=1071\", \"36\u00d713=468\"],\n  [\"62\u00d741=2542\", \"84\u00d752=4368\"],\n  [\"59\u00d725=1475\", \"94\u00d712=1128\"],\n  [\"93\u00d728=2604\", \"86\u00d737=3182\"],\n  [\"30\u00d766=1980\", \"97\u00d750=4850\"],\n  [\"52\u00d713=676\", \"41\u00d770=2870\"],\n  [\"27\u00d764=1728\", \"79\u00d771=5609\"],\n  [\"64\u00d735=2240\", \"97\u00d769=6693\"],\n  [\"13\u00d716=208\", \"41\u00d762=2542\"],\n  [\"53\u00d779=4187\", \"89\u00d766=5874\"],\n  [\"55\u00d783=4565\", \"46\u00d727=1242\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date heading and every two-digit-by-two-digit multiplication\n# problem/answer in the table with the new values from the target revision.\n# Every \"old\" value below is unique in the document, so a single\n# Find/Replace (wrap = none, replace one occurrence) per pair is\n# sufficient and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-28 Wednesday\", \"2025-05-29 Thursday\"),\n    @(\"64\u00d784=5376\", \"65\u00d747=3055\"),\n    @(\"63\u00d758=3654\", \"45\u00d789=4005\"),\n    @(\"91\u00d723=2093\", \"36\u00d796=3456\"),\n    @(\"14\u00d755=770\", \"20\u00d728=560\"),\n    @(\"63\u00d738=2394\", \"45\u00d763=2835\"),\n    @(\"18\u00d713=234\", \"18\u00d773=1314\"),\n    @(\"96\u00d782=7872\", \"59\u00d781=4779\"),\n    @(\"43\u00d714=602\", \"61\u00d755=3355\"),\n    @(\"46\u00d731=1426\", \"71\u00d745=3195\"),\n    @(\"93\u00d747=4371\", \"71\u00d779=5609\"),\n    @(\"72\u00d798=7056\", \"13\u00d718=234\"),\n    @(\"97\u00d788=8536\", \"37\u00d711=407\"),\n    @(\"78\u00d768=5304\", \"17\u00d720=340\"),\n    @(\"76\u00d747=3572\", \"28\u00d793=2604\"),\n    @(\"17\u00d763=1071\", \"36\u00d713=468\"),\n    @(\"62\u00d741=2542\", \"84\u00d752=4368\"),\n    @(\"59\u00d725=1475\", \"94\u00d712=1128\"),\n    @(\"93\u00d728=2604\", \"86\u00d737=3182\"),\n    @(\"30\u00d766=1980\", \"97\u00d750=4850\"),\n    @(\"52\u00d713=676\", \"41\u00d770=2870\"),\n    @(\"27\u00d764=1728\", \"79\u00d771=5609\"),\n    @(\"64\u00d735=2240\", \"97\u00d769=6693\"),\n    @(\"13\u00d716=208\", \"41\u00d762=2542\"),\n    @(\"53\u00d779=4187\", \"89\u00d766=5874\"),\n    @(\"55\u00d783=4565\", \"46\u00d727=1242\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, $false, $true, $false, $false, $false, $true, 1, $false, [ref]$newText, 2)\n}\n"}
